# fix: add user setting sample APi seqs
# Shift rows 2-10 down by one logical "record" (values move up the table by
# one row each, i.e. row N takes on the values that used to belong further
# down), refresh the timestamp on every existing data row, and append a new
# row 11 with the data that used to live in the old row 10 (the DELETE
# /api/v1/memo/21 record), extending the sheet's used range to A1:Q11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-06 12:03:12"

# Data for rows 2..11 (A..Q), in final desired state.
$rows = @(
    @{ Row=2;  A=5;  B=4;  C=5;  D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/resource/16"; G="/api/v1/resource/16"; H="{}"; I=""; J=2; K=5; L=200; M=0.003; N=0; O=0; P=$true;  Q=$true }
    @{ Row=3;  A=6;  B=5;  C=6;  D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/memo/21";     G="/api/v1/memo/21";     H="{}"; I=""; J=2; K=5; L=200; M=0.002; N=1; O=0; P=$true;  Q=$false }
    @{ Row=4;  A=7;  B=5;  C=7;  D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/resource/16"; G="/api/v1/resource/16"; H="{}"; I=""; J=2; K=5; L=200; M=0.002; N=1; O=1; P=$false; Q=$false }
    @{ Row=5;  A=8;  B=6;  C=8;  D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/resource/16"; G="/api/v1/resource/16"; H="{}"; I=""; J=2; K=5; L=200; M=0.002; N=0; O=0; P=$true;  Q=$true }
    @{ Row=6;  A=9;  B=7;  C=9;  D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/memo/21";     G="/api/v1/memo/21";     H="{}"; I=""; J=2; K=5; L=200; M=0.002; N=1; O=0; P=$true;  Q=$false }
    @{ Row=7;  A=10; B=7;  C=10; D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/resource/16"; G="/api/v1/resource/16"; H="{}"; I=""; J=2; K=5; L=200; M=0.002; N=1; O=1; P=$false; Q=$false }
    @{ Row=8;  A=11; B=8;  C=11; D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/resource/16"; G="/api/v1/resource/16"; H="{}"; I=""; J=2; K=5; L=200; M=0.003; N=0; O=0; P=$true;  Q=$true }
    @{ Row=9;  A=12; B=9;  C=12; D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/resource/16"; G="/api/v1/resource/16"; H="{}"; I=""; J=2; K=5; L=200; M=0.003; N=1; O=1; P=$false; Q=$false }
    @{ Row=10; A=13; B=10; C=13; D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/memo/21";     G="/api/v1/memo/21";     H="{}"; I=""; J=2; K=5; L=200; M=0.003; N=0; O=0; P=$true;  Q=$true }
    @{ Row=11; A=14; B=11; C=14; D=$newTimestamp; E="DELETE"; F="http://49.234.6.241:5230/api/v1/memo/21";     G="/api/v1/memo/21";     H="{}"; I=""; J=2; K=5; L=200; M=0.002; N=1; O=1; P=$false; Q=$false }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    # Column I ("data") is an empty inline string on every existing data row
    # (r2..r10) already, so it is intentionally left untouched here - writing
    # "" through COM drops the cell instead of keeping it blank. Row 11 is
    # brand new, so its blank "data" cell is created explicitly below.
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
}

# Row 11 is new, so I11 doesn't exist yet as the blank placeholder text cell
# that every other data row has (I2:I10). Writing "" directly clears/omits
# the cell instead of keeping it as blank text, so force a text-typed empty
# value the way Excel's UI does (leading apostrophe -> empty text, not a
# literal character in the stored value), then drop the quote-prefix style
# that trick leaves behind so the cell keeps the default (unstyled) look.
$ws.Cells.Item(11, 9).Value = "'"
$ws.Cells.Item(11, 9).Style = "Normal"
